# TestDataLazada.xlsx - "test login - done!"
#
# The login-form test sheet ("2.2. DN_Form dang nhap") had its test-case
# results filled in/corrected, and the "linked account" sheet
# ("2.3. DN_Tai khoan lien ket") had its previously-blank Status column
# filled in with results. The active sheet/selection in the workbook was
# also left on the login-form sheet, matching where the tester was last
# working.

$wb = $excel.ActiveWorkbook

$wsLinks   = $wb.Worksheets.Item(2)   # "2.1. DN_Kiem tra links"
$wsLogin   = $wb.Worksheets.Item(3)   # "2.2. DN_Form dang nhap"
$wsLinked  = $wb.Worksheets.Item(4)   # "2.3. DN_Tai khoan lien ket"

# --- "2.2. DN_Form dang nhap": test case DN_14 (row 8) actually failed ---
$wsLogin.Activate()
$wsLogin.Range("G8").Value = "Fail"
$wsLogin.Range("E7:E8").Select()

# --- "2.1. DN_Kiem tra links": selection left on C2, no longer the active tab ---
$wsLinks.Range("C2").Select()

# --- "2.3. DN_Tai khoan lien ket": fill in the Status column results ---
$wsLinked.Range("F2").Value = "Pass"
$wsLinked.Range("F3").Value = "Pass"
$wsLinked.Range("C2").Select()

# Leave the login-form sheet as the active/selected tab, as in the committed workbook
$wsLogin.Activate()
$wsLogin.Range("E7:E8").Select()
